# Apply the FHIR IG spreadsheet regeneration changes:
#  - Metadata sheet: bump Version, Date, set Publisher, replace the duplicated
#    "Contact" row with a "Jurisdiction" row, and drop the trailing duplicate
#    row so the sheet shrinks from 21 to 20 rows.
#  - Elements sheet: the root Extension's Short/Definition cells get the
#    real title/description instead of the generic "Extension"/"An Extension"
#    placeholders.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Remove the second, duplicated "Contact" / "No display for ContactDetail" row
# (old row 11). Everything below shifts up by one, shrinking the sheet from
# 21 to 20 rows, matching the new dimension A1:B20.
$meta.Rows.Item(11).Delete()

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: regenerated timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a value
$meta.Range("B9").Value = "Alvearie Team"

# The remaining "Contact" row (old row 10, now still row 10 after the delete)
# becomes the new "Jurisdiction" row
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

$elements = $wb.Worksheets.Item("Elements")

# Root Extension row: Short / Definition get the real title & description
$elements.Range("K2").Value = "Capitated Service Indicator"
$elements.Range("L2").Value = "Indicates whether this service (encounter record) was capitated"
